$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G: new "meeting" date column (03 Jan 2024) ---
# Give G1 the same date-number-format as the other date header cells (C1:F1)
$ws.Range("G1").NumberFormat = $ws.Range("F1").NumberFormat
$ws.Range("G1").Value = 45294

# Attendance for the existing students (rows 2-8) on the new date.
# Everyone is Present except row 7 (Minakshi Kadao).
$ws.Range("G2").Value = "Present"
$ws.Range("G3").Value = "Present"
$ws.Range("G4").Value = "Present"
$ws.Range("G5").Value = "Present"
$ws.Range("G6").Value = "Present"

# --- New row 9: a new student, Ayesha Singh (entered before marking row 7
# absent so the shared-string table gains "Ayesha Singh" before "Absent") ---
# Copy formatting (styles) from row 8 down to row 9 first.
$ws.Range("A8:T8").Copy()
$ws.Range("A9:T9").PasteSpecial(-4122)
$ws.Rows.Item(9).RowHeight = $ws.Rows.Item(8).RowHeight

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Ayesha Singh"
$ws.Range("C9").Value = "Absent"
$ws.Range("D9").Value = "Absent"
$ws.Range("E9").Value = "Absent"
$ws.Range("F9").Value = "Absent"
$ws.Range("G9").Value = "Present"

# Row 7 (Minakshi Kadao) is marked Absent for the new date.
$ws.Range("G7").Value = "Absent"
$ws.Range("G8").Value = "Present"

# --- Extend the attendance dropdown validation down to the new row ---
$ws.Range("C2:N8").Validation.Delete()
$ws.Range("C2:N9").Validation.Add(3, 1, 1, '"Present, Absent,Reason"')

Write-Output "done"
